# Directorio.xlsx update
# - Update the description in F25 (Bodega de Tiempos entry)
# - Insert two new rows (26, 27) describing "Programaciones por capitulo" and
#   "Modelo de Datos"
# - Add a hyperlink for the new Colab URL in G26
# - Grow the "Tabla13" Excel table to cover the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the existing "Bodega de Tiempos" row (row 25): the description in
#    F25 now documents the in-progress Colab cleanup task.
# ---------------------------------------------------------------------------
$ws.Range("F25").Value = "En proceso: Colab limpieza de Bodega de tiempos con cargue desde local de la carpeta de ETL y integración con nombres de proyectos"

# ---------------------------------------------------------------------------
# 2. New row 26: "5.1.2. Programaciones por capitulo" (child of 5.1.)
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = "5.1.2."
$ws.Range("B26").Value = "5.1."
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = "Programaciones por capitulo"
$ws.Range("E26").Value = "url"
$ws.Range("F26").Value = "Finalizado: Colab de vetificación de incidencia en cierre de capitulos de workplanner ( se revisar por programación, se puede ampliar a bodega de tiempos)"
$ws.Range("G26").Value = "https://colab.research.google.com/drive/1XhnifKohqfhBAwTACpzS7d6D4qcx_sl0#scrollTo=3ZOFGD83yd0f"

# Hyperlink for the new Colab URL (split into address + sub-address, matching
# the existing workbook convention for "#scrollTo=" style links).
$ws.Hyperlinks.Add($ws.Range("G26"), "https://colab.research.google.com/drive/1XhnifKohqfhBAwTACpzS7d6D4qcx_sl0", "scrollTo=3ZOFGD83yd0f")

# ---------------------------------------------------------------------------
# 3. New row 27: "5.2. Modelo de Datos" (child of 5.)
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "5.2."
# "5." alone would be auto-coerced to the number 5 by the usual Value
# assignment (Excel's numeric-text inference), so copy it as text from the
# existing "5." cell (B24) instead of typing it.
$ws.Range("B24").Copy()
$ws.Range("B27").PasteSpecial(-4163)
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = "Modelo de Datos"
$ws.Range("E27").Value = "Carpeta"
$ws.Range("F27").Value = "Se plantea modelo de datos bajo trabajo en Bigquery y Colab"

# ---------------------------------------------------------------------------
# 4. Formatting: copy the look of similar existing rows onto the new ones.
#    Row 20 (A20:H20) has the same shape as row 26 (data row with a
#    hyperlinked URL cell); row 24 (A24:H24) has the same shape as row 27
#    (data row with no URL/hyperlink).
# ---------------------------------------------------------------------------
$ws.Range("A20:H20").Copy()
$ws.Range("A26:H26").PasteSpecial(-4122)

$ws.Range("A24:H24").Copy()
$ws.Range("A27:H27").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. Grow the worksheet table ("Tabla13") so it covers the new rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H27"))

# ---------------------------------------------------------------------------
# 6. Update the active selection to reflect where editing ended (D27), as in
#    the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("D27").Select()
